$d = $word.ActiveDocument

# --- Step 1: merge run pairs for TaxesWoCredit2 / Taxes2 / EnergyRevenueRequired2 ---
$d.Content.Find.Execute("TaxesWoCredit2 =((CombinedTaxRate", $true, $false, $false, $false, $false, $true, 1, $false, "TaxesWoCredit2 =((CombinedTaxRate", 2) | Out-Null
$d.Content.Find.Execute("Taxes2 =((CombinedTaxRate", $true, $false, $false, $false, $false, $true, 1, $false, "Taxes2 =((CombinedTaxRate", 2) | Out-Null
$d.Content.Find.Execute("EnergyRevenueRequired2 =EquityRecovery", $true, $false, $false, $false, $false, $true, 1, $false, "EnergyRevenueRequired2 =EquityRecovery", 2) | Out-Null

# --- Step 2: remove the _GoBack bookmark paragraph (keep it as a blank paragraph) ---
$bookmarkPara = $d.Paragraphs.Item(105)
$bookmarkRange = $bookmarkPara.Range
$blankXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$bookmarkRange.InsertXML($blankXml)

# --- Step 3: append the new "Current $ Level Annual Cost (LAC)" section at the end ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newSectionXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
      <w:pPr>
        <w:pStyle w:val="Heading3"/>
      </w:pPr>
      <w:r>
        <w:t>Current $ Level Annual Cost (LAC)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>PresentWorth1=EnergyRevenueRequired1*(1+CostOfEquity/</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>100)^</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>-B130</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TotalPresentWorth</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>SUM(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>PresentWorth1</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>PresentWorth</w:t>
      </w:r>
      <w:r>
        <w:t>20)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactorCurrent</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactor</w:t>
      </w:r>
      <w:r>
        <w:t>Equity</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CurrentLevelAnnualRevenueRequirements</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TotalPresentWorth</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>*</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactorCurrent</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CurrentLACofEnergy</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CurrentLevelAnnualRevenueRequirements</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>AnnualGeneration</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>RealCostOfMoney</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>=(1+CostOfEquity/100)/(1+GeneralInflation/100)-1</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>CapitalRecoveryFactorConstant</w:t>
      </w:r>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactor</w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>RealCostOfMone</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">y*100, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>EconomicLife</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ConstantLevelAnnualRevenueRequirements</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TotalPresentWorth</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>*</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactorConstant</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ConstantLACofEnergy</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ConstantLevelAnnualRevenueRequirements</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>AnnualGeneration</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    function </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>CapitalRecoveryFactor</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: number, N: number) {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">        const A = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> / 100 * </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Math.pow</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">(1 + </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> / 100), N) / (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Math.pow</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">((1 + </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> / 100), N) - 1);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">        return A;</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$endRange.InsertXML($newSectionXml)
